# Recalibration update: "Settings back to qfin22 presentation - recalibration"
# Updates the Linear and NonLinear parameter sheets with new calibration values.

$wb = $excel.ActiveWorkbook

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")

$wsLinear.Range("B2").Value = 0.02820800260415238
$wsLinear.Range("B3").Value = -0.01537376613605469
$wsLinear.Range("B4").Value = 47.6528387904176
$wsLinear.Range("B5").Value = "[0.9999999999999999, 0.20260799747081493, 0.04919150166704316, 0.05429277802523905, 0.04736967659941702, 0.03798439392662625, 0.2265122473035163, 0.3274653924481518, 0.1859990769625859, 0.03820672099433051, -0.005769301243176278, 0.032318980887897045, 0.04734943962213661, 0.1830632214138856, 0.3240974299002695, 0.19759061168874995, 0.029844494167913415, 0.017033788417185784, 0.01892646979476015, 0.019780234049373163]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")

$wsNonLinear.Range("B4").Value = -0.3084452739430301
$wsNonLinear.Range("B5").Value = -0.0839750369937859
$wsNonLinear.Range("B6").Value = 48.84702065543422
$wsNonLinear.Range("B7").Value = 0.408346165694867
$wsNonLinear.Range("B8").Value = -0.1138110822101487
$wsNonLinear.Range("B9").Value = 46.43854841019247
$wsNonLinear.Range("B10").Value = "[1.0, 0.19880933797324749, 0.05248618169215941, 0.05728935354265476, 0.05182603610508294, 0.04286620081499489, 0.2248461022484913, 0.32026542189281054, 0.18410791935439375, 0.0428514506702709, -0.0009480066279570798, 0.03634546821424419, 0.05242379666624795, 0.18128553563518357, 0.3163383837106515, 0.19568527533802135, 0.03338687213494882, 0.021865145105055328, 0.023718508008376923, 0.0238369612856056]"
